# Add a new "GARCH" column (E) to the comparison table, mirroring the
# style of the existing header row and filling in the metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching style of other header cells (B1:D1)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value = "GARCH"

# Data values for rows 2-17
$ws.Range("E2").Value = 0.43
$ws.Range("E3").Value = 0.65
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 1219584.15
$ws.Range("E6").Value = 2408660.67
$ws.Range("E7").Value = 1.975
$ws.Range("E8").Value = 0.1821
$ws.Range("E9").Value = 0.2944
$ws.Range("E10").Value = 1.77
$ws.Range("E11").Value = 3
$ws.Range("E12").Value = 4549
$ws.Range("E13").Value = 0.0004
$ws.Range("E14").Value = 0.7610343290236291
$ws.Range("E15").Value = 0.1489108761279314
$ws.Range("E16").Value = -0.06610456994113106
$ws.Range("E17").Value = 0.1122050666664724
